$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") - reuse the same header
# formatting (bold font, border, centered/top alignment) already applied
# to B1:H1 by copying formats only from H1.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data columns I (I0) and J (IF) for rows 2-12.
$values = @{
    2  = @(6, 7)
    3  = @(8, 9)
    4  = @(7, 9)
    5  = @(10, 10)
    6  = @(7, 7)
    7  = @(5, 7)
    8  = @(8, 8)
    9  = @(8, 9)
    10 = @(1, 3)
    11 = @(8, 9)
    12 = @(1, 1)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
